$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number (single decimal point) need
# to be pinned to Text format first, otherwise Excel auto-converts the literal
# into a numeric value (dropping trailing zeros / introducing float noise).
$textCells = @("D5", "D6", "D7", "D8", "D9", "D10", "D12", "D14", "D18", "D20", "D21", "D22", "D23", "D24", "D26", "D27", "D29", "D30", "D32", "D34", "D35", "D36", "D38", "D39", "D40", "D41", "D42", "D43", "D45", "D47", "D49", "D50", "D51")
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = "96.599.99"
$ws.Range("E2").Value = "  -0.05%  "
$ws.Range("D3").Value = "3.629.78"
$ws.Range("E3").Value = "  +1.12%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "242.51"
$ws.Range("E5").Value = "  +0.18%  "
$ws.Range("D6").Value = "1.80"
$ws.Range("E6").Value = "  +15.36%  "
$ws.Range("D7").Value = "655.32"
$ws.Range("E7").Value = "  -0.53%  "
$ws.Range("D8").Value = "0.421"
$ws.Range("E8").Value = "  +3.32%  "
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").Value = "1.07"
$ws.Range("E9").Value = "  +1.33%  "
$ws.Range("B10").Value = "USDC"
$ws.Range("C10").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D10").Value = "1.00"
$ws.Range("E10").Value = "  -0.04%  "
$ws.Range("D11").Value = "3.619.45"
$ws.Range("E11").Value = "  +0.87%  "
$ws.Range("D12").Value = "44.30"
$ws.Range("E12").Value = "  +2.15%  "
$ws.Range("E13").Value = "  +0.41%  "
$ws.Range("D14").Value = "6.51"
$ws.Range("E14").Value = "  +1.74%  "
$ws.Range("D15").Value = "4.302.69"
$ws.Range("E15").Value = "  +1.04%  "
$ws.Range("D16").Value = "96.376.68"
$ws.Range("E16").Value = "  -0.19%  "
$ws.Range("E17").Value = "  -0.24%  "
$ws.Range("D18").Value = "8.64"
$ws.Range("E18").Value = "  +11.32%  "
$ws.Range("D19").Value = "3.634.38"
$ws.Range("E19").Value = "  +1.19%  "
$ws.Range("D20").Value = "13.09"
$ws.Range("E20").Value = "  +3.65%  "
$ws.Range("D21").Value = "18.41"
$ws.Range("E21").Value = "  +3.01%  "
$ws.Range("D22").Value = "0.527"
$ws.Range("E22").Value = "  +6.51%  "
$ws.Range("B23").Value = "SuiNetwork"
$ws.Range("C23").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D23").Value = "3.44"
$ws.Range("E23").Value = "  -0.41%  "
$ws.Range("B24").Value = "BitcoinCash"
$ws.Range("C24").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D24").Value = "511.36"
$ws.Range("E24").Value = "  -0.41%  "
$ws.Range("D26").Value = "6.88"
$ws.Range("E26").Value = "  +0.58%  "
$ws.Range("D27").Value = "100.95"
$ws.Range("E27").Value = "  +4.27%  "
$ws.Range("E28").Value = "  +3.01%  "
$ws.Range("D29").Value = "0.166"
$ws.Range("E29").Value = "  +11.02%  "
$ws.Range("D30").Value = "3.04"
$ws.Range("E30").Value = "  +0.49%  "
$ws.Range("E31").Value = "  +2.99%  "
$ws.Range("D32").Value = "0.999"
$ws.Range("E32").Value = "  -0.11%  "
$ws.Range("E33").Value = "  -0.69%  "
$ws.Range("D34").Value = "33.11"
$ws.Range("E34").Value = "  +4.33%  "
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.51%  "
$ws.Range("D36").Value = "1.73"
$ws.Range("E36").Value = "  +7.79%  "
$ws.Range("E37").Value = "  +1.97%  "
$ws.Range("D38").Value = "8.85"
$ws.Range("E38").Value = "  +4.04%  "
$ws.Range("D39").Value = "616.39"
$ws.Range("E39").Value = "  +2.78%  "
$ws.Range("D40").Value = "41.52"
$ws.Range("E40").Value = "  +20.57%  "
$ws.Range("D41").Value = "0.156"
$ws.Range("D42").Value = "0.952"
$ws.Range("E42").Value = "  +4.55%  "
$ws.Range("D43").Value = "1.94"
$ws.Range("E43").Value = "  +4.99%  "
$ws.Range("D45").Value = "6.13"
$ws.Range("E45").Value = "  +6.50%  "
$ws.Range("E46").Value = "  +3.91%  "
$ws.Range("D47").Value = "2.30"
$ws.Range("E47").Value = "  +0.80%  "
$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").Value = "8.54"
$ws.Range("E49").Value = "  +3.20%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "0.401"
$ws.Range("E50").Value = "  +14.02%  "
$ws.Range("D51").Value = "54.61"
$ws.Range("E51").Value = "  +1.57%  "
